# Applies the "last changes model characteristics" commit to the workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("model_characteristics")
$ws2 = $wb.Worksheets.Item("model_overview")

# --- model_characteristics (sheet1) data updates -----------------------

# Observables (col C) corrections
$ws1.Range("C2").Value  = 11
$ws1.Range("C13").Value = 3

# Conditions (col E) correction
$ws1.Range("E12").Value = 17

# Steady State (col L) values incremented
$ws1.Range("L2").Value  = 1
$ws1.Range("L3").Value  = 2
$ws1.Range("L4").Value  = 1
$ws1.Range("L5").Value  = 2
$ws1.Range("L6").Value  = 3
$ws1.Range("L7").Value  = 1
$ws1.Range("L10").Value = 2
$ws1.Range("L11").Value = 1
$ws1.Range("L12").Value = 2
$ws1.Range("L13").Value = 3
$ws1.Range("L14").Value = 1
$ws1.Range("L15").Value = 2
$ws1.Range("L16").Value = 1
$ws1.Range("L17").Value = 2
$ws1.Range("L18").Value = 3
$ws1.Range("L19").Value = 2
$ws1.Range("L20").Value = 3
$ws1.Range("L21").Value = 3

# Legend for "Error Model" (rows 25-28, column H) shifts one position,
# replacing the obsolete "Fixed errors" entry with a new "Ex from paper"
# one at the end. Set H28 first so the newly created shared strings are
# appended in the same order as in the target workbook.
$ws1.Range("H28").Value = "E3 from paper"
$ws1.Range("H25").Value = "Ex from paper"
$ws1.Range("H26").Value = "E1 from paper"
$ws1.Range("H27").Value = "E2 from paper"

# --- model_overview (sheet2) mirrored data updates ----------------------

$ws2.Range("D3").Value  = 11
$ws2.Range("E13").Value = 17
$ws2.Range("D14").Value = 3

# --- active sheet / selection bookkeeping -------------------------------

[void]$ws1.Range("C3").Select()
$ws1.Activate()

[void]$ws2.Range("D4").Select()
$ws2.Activate()
